$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" '26.707.87'
$ws.Range("E2").Value = '  +0.28%  '
Set-TextValue "D3" '1.598.22'
$ws.Range("E3").Value = '  +0.09%  '
$ws.Range("E4").Value = '  +0.19%  '
Set-TextValue "D5" '211.10'
$ws.Range("E5").Value = '  -0.05%  '
$ws.Range("E6").Value = '  -0.70%  '
$ws.Range("E7").Value = '  +0.23%  '
$ws.Range("E8").Value = '  +0.28%  '
$ws.Range("E9").Value = '  +0.94%  '
$ws.Range("E10").Value = '  +0.74%  '
$ws.Range("E11").Value = '  +0.43%  '
Set-TextValue "D12" '1.823.38'
$ws.Range("E12").Value = '  +0.19%  '
Set-TextValue "D13" '1.606.66'
$ws.Range("E13").Value = '  +0.74%  '
$ws.Range("E14").Value = '  +0.46%  '
$ws.Range("E15").Value = '  +0.27%  '
Set-TextValue "D16" '65.27'
$ws.Range("E16").Value = '  +0.86%  '
Set-TextValue "D17" '26.684.95'
$ws.Range("E17").Value = '  +0.29%  '
Set-TextValue "D18" '0.0₃0761'
$ws.Range("E18").Value = '  +4.17%  '
$ws.Range("E19").Value = '  +0.96%  '
$ws.Range("E20").Value = '  +0.21%  '
Set-TextValue "D21" '7.12'
$ws.Range("E21").Value = '  +2.71%  '
Set-TextValue "D22" '4.27'
$ws.Range("E22").Value = '  +0.33%  '
$ws.Range("E23").Value = '  -0.57%  '
Set-TextValue "D24" '8.92'
$ws.Range("E24").Value = '  +0.70%  '
Set-TextValue "D25" '143.16'
$ws.Range("E25").Value = '  -1.59%  '
$ws.Range("E26").Value = '  +0.19%  '
$ws.Range("E27").Value = '  -0.31%  '
$ws.Range("E28").Value = '  +0.21%  '
Set-TextValue "D29" '15.31'
$ws.Range("E29").Value = '  +0.05%  '
Set-TextValue "D30" '0.0517'
$ws.Range("E30").Value = '  +2.68%  '
$ws.Range("E31").Value = '  -0.27%  '
$ws.Range("E32").Value = '  +0.32%  '
$ws.Range("E33").Value = '  +1.52%  '
Set-TextValue "D34" '1.290.11'
$ws.Range("E34").Value = '  +0.51%  '
$ws.Range("E35").Value = '  -5.82%  '
$ws.Range("E36").Value = '  +0.97%  '
Set-TextValue "D37" '1.49'
$ws.Range("E37").Value = '  +0.11%  '
$ws.Range("E38").Value = '  -0.26%  '
$ws.Range("E39").Value = '  +17.66%  '
Set-TextValue "D40" '0.825'
$ws.Range("E40").Value = '  -1.73%  '
$ws.Range("E41").Value = '  -0.03%  '
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue "D42" '0.784'
$ws.Range("E42").Value = '  -0.29%  '
$ws.Range("B43").Value = 'MXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue "D43" '2.19'
$ws.Range("E43").Value = '  -0.34%  '
Set-TextValue "D44" '63.02'
$ws.Range("E44").Value = '  -1.12%  '
Set-TextValue "D45" '1.734.58'
$ws.Range("E45").Value = '  +0.13%  '
Set-TextValue "D46" '91.23'
$ws.Range("E46").Value = '  +1.69%  '
$ws.Range("E47").Value = '  -1.53%  '
$ws.Range("E48").Value = '  -0.95%  '
$ws.Range("E49").Value = '  +0.57%  '
$ws.Range("E50").Value = '  +0.12%  '
Set-TextValue "D51" '7.33'
$ws.Range("E51").Value = '  -1.69%  '
